# Applies the "Render DB env hardening" update to the requirement-status doc.
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.ClearFormatting()
    $ok = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $old"
    }
}

# --- 1) In-place rewrites of the first 5 paragraphs of the "New in this update" block ---
Replace-Text "New in this update (Render React peer fix)" "New in this update (Render DB env hardening)"

Replace-Text "- Fixed frontend dependency conflict causing Render build failure:" "- Updated backend DB connection resolver in ``Program.cs`` to support multiple env keys:"

Replace-Text "  - ``react-day-picker@8.10.1`` supports React up to v18." "  - ``POSTGRES_CONNECTION_STRING``"

Replace-Text "  - Downgraded ``react`` and ``react-dom`` from ``^19.0.0`` to ``^18.2.0`` in ``frontend/package.json``." "  - ``ConnectionStrings__Postgres`` / ``ConnectionStrings:Postgres``"

Replace-Text "- This resolves ERESOLVE error for ``react-day-picker`` peer dependency on Render." "  - ``DATABASE_URL``"

# --- 2) Insert the 4 brand-new paragraphs after the (now renamed) DATABASE_URL bullet ---
# The 5 rewrites above are 1:1 text swaps, so paragraph count/order is unchanged;
# the DATABASE_URL bullet is still the paragraph that used to read the ERESOLVE line
# (the 5th paragraph of the "New in this update" block, i.e. document paragraph 16).
$anchorPara = $d.Paragraphs(16)
if ($anchorPara.Range.Text.TrimEnd([char]13) -ne "  - ``DATABASE_URL``") {
    throw "Anchor paragraph 16 did not contain the expected DATABASE_URL bullet"
}

$newLines = @(
    "  - ``RENDER_EXTERNAL_DATABASE_URL``",
    "  - ``RENDER_INTERNAL_DATABASE_URL``",
    "- Added support to normalize ``postgres://...`` URLs into Npgsql connection string format.",
    "- Added production safeguard: if resolved host is localhost, app fails fast with clear error message."
)

$insertIndex = 16
foreach ($line in $newLines) {
    $d.Paragraphs($insertIndex).Range.InsertParagraphAfter() | Out-Null
    $insertIndex = $insertIndex + 1
    $d.Paragraphs($insertIndex).Range.Text = $line
}

# --- 3) Git state section updates ---
Replace-Text "- Last pushed commit: 7d0fd23" "- Last pushed commit: 0c88abb"

Replace-Text "- Current React compatibility fix is local and not pushed yet." "- Current DB env resilience fix is local and not pushed yet."

Write-Output "edit complete"
